# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
# for the coinranking.com snapshot on Sheet1.
#
# Column D cells are stored as plain text (e.g. "61.983.83", "1.00") even
# when they look like numbers, so values that parse as a plain number are
# written back through a temporary "Text" number format (restoring the
# cell's original style afterwards) to stop Excel from re-interpreting them
# as numeric values and dropping significant trailing/leading zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '61.983.83'; E = '  +2.79%  ' },
    @{ Row = 3; D = '3.412.67'; E = '  +3.48%  ' },
    @{ Row = 4; D = '1.00'; E = '  -0.06%  ' },
    @{ Row = 5; D = '577.54'; E = '  +2.93%  ' },
    @{ Row = 6; D = '138.41'; E = '  +7.38%  ' },
    @{ Row = 7; E = '  -0.09%  ' },
    @{ Row = 8; D = '3.412.47'; E = '  +3.48%  ' },
    @{ Row = 9; E = '  +1.25%  ' },
    @{ Row = 10; D = '7.52'; E = '  +2.30%  ' },
    @{ Row = 11; D = '0.128'; E = '  +9.55%  ' },
    @{ Row = 12; E = '  +6.99%  ' },
    @{ Row = 13; D = '3.992.97'; E = '  +3.19%  ' },
    @{ Row = 14; E = '  +1.85%  ' },
    @{ Row = 15; E = '  +8.27%  ' },
    @{ Row = 16; D = '3.414.49'; E = '  +3.46%  ' },
    @{ Row = 17; D = '25.50'; E = '  +5.52%  ' },
    @{ Row = 18; D = '61.994.69'; E = '  +2.39%  ' },
    @{ Row = 19; E = '  +6.11%  ' },
    @{ Row = 20; E = '  +4.98%  ' },
    @{ Row = 21; D = '9.49'; E = '  +6.12%  ' },
    @{ Row = 22; D = '390.57'; E = '  +11.37%  ' },
    @{ Row = 23; E = '  +3.65%  ' },
    @{ Row = 24; D = '3.548.61'; E = '  +3.32%  ' },
    @{ Row = 25; E = '  +19.13%  ' },
    @{ Row = 26; E = '  +0.06%  ' },
    @{ Row = 27; D = '71.58'; E = '  +3.41%  ' },
    @{ Row = 28; D = '1.61'; E = '  +11.48%  ' },
    @{ Row = 29; D = '7.65'; E = '  +4.98%  ' },
    @{ Row = 30; D = '0.996'; E = '  -0.30%  ' },
    @{ Row = 31; D = '8.33'; E = '  +6.42%  ' },
    @{ Row = 32; E = '  +5.37%  ' },
    @{ Row = 33; E = '  +3.39%  ' },
    @{ Row = 34; D = '3.442.43'; E = '  +3.33%  ' },
    @{ Row = 35; E = '  -0.02%  ' },
    @{ Row = 36; E = '  +4.16%  ' },
    @{ Row = 37; D = '5.55'; E = '  +5.34%  ' },
    @{ Row = 38; D = '6.99'; E = '  +3.40%  ' },
    @{ Row = 39; E = '  +5.15%  ' },
    @{ Row = 40; D = '161.49'; E = '  +2.29%  ' },
    @{ Row = 41; D = '0.0796'; E = '  +6.06%  ' },
    @{ Row = 42; D = '1.75'; E = '  +13.79%  ' },
    @{ Row = 43; D = '1.00'; E = '  -0.09%  ' },
    @{ Row = 44; E = '  +6.75%  ' },
    @{ Row = 45; D = '0.777'; E = '  +4.89%  ' },
    @{ Row = 46; E = '  +3.28%  ' },
    @{ Row = 47; D = '25.29'; E = '  +11.72%  ' },
    @{ Row = 48; D = '41.61'; E = '  +1.80%  ' },
    @{ Row = 49; D = '6.99'; E = '  +4.71%  ' },
    @{ Row = 50; D = '22.95'; E = '  +6.39%  ' },
    @{ Row = 51; D = '2.391.22'; E = '  +10.47%  ' }
)

foreach ($item in $updates) {
    $row = $item.Row

    if ($item.ContainsKey("D")) {
        $newValue = $item.D
        $cell = $ws.Range("D$row")

        if ($newValue -match '^-?[0-9]+(\.[0-9]+)?$') {
            # Looks like a plain number (e.g. "1.00", "7.52") - force text
            # storage so Excel keeps it verbatim instead of coercing it to
            # a numeric value, then restore the original (unstyled) look.
            $originalStyle = $cell.Style
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
            $cell.Style = $originalStyle
        } else {
            $cell.Value = $newValue
        }
    }

    if ($item.ContainsKey("E")) {
        $ws.Range("E$row").Value = $item.E
    }
}
